$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header: "Municipio" -> "CVE_MUN"
$ws.Range("A1").Value = "CVE_MUN"

# Municipality names in A2:A85 -> INEGI municipality keys 13001..13084.
# Format the range as text first so the numeric-looking keys are stored
# as literal text (matching the source data), not auto-converted numbers.
$codeRange = $ws.Range("A2:A85")
$codeRange.NumberFormat = "@"

$values = New-Object 'object[,]' 84,1
for ($i = 0; $i -lt 84; $i++) {
    $values[$i, 0] = [string](13001 + $i)
}
$codeRange.Value = $values
